$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4234
$ws.Range("I64").Value = 5103.846
$ws.Range("J64").Value = 3291.6667
$ws.Range("K64").Value = 5103.846
$ws.Range("L64").Value = 3291.6667
$ws.Range("M64").Value = -4855.846
$ws.Range("N64").Value = -3787.6667

$ws.Range("H67").Value = 4234
$ws.Range("I67").Value = 5103.846
$ws.Range("J67").Value = 3291.6667
$ws.Range("K67").Value = 5103.846
$ws.Range("L67").Value = 3291.6667
$ws.Range("M67").Value = -4245.846
$ws.Range("N67").Value = -5007.6667

$ws.Range("H74").Value = 5437.4546
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064

$ws.Range("H77").Value = 5437.4546
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320

$ws.Range("H113").Value = 31252252
$ws.Range("J113").Value = 41668668
$ws.Range("L113").Value = 41668668
$ws.Range("N113").Value = -41675176

$ws.Range("H138").Value = 3029.8452
$ws.Range("I138").Value = 1492
$ws.Range("J138").Value = 4183.229
$ws.Range("K138").Value = 4476
$ws.Range("L138").Value = 12549.687
$ws.Range("M138").Value = 664
$ws.Range("N138").Value = -22829.687

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16428.955
$ws.Range("I32").Value = 12386.667
$ws.Range("J32").Value = 22849.059
$ws.Range("K32").Value = 12386.667
$ws.Range("L32").Value = 22849.059
$ws.Range("M32").Value = -12099.667
$ws.Range("N32").Value = -23423.059

$ws.Range("H63").Value = 166677630
$ws.Range("I63").Value = 333335260
$ws.Range("K63").Value = 333335260
$ws.Range("M63").Value = -333334574

$ws.Range("H66").Value = 166677630
$ws.Range("I66").Value = 333335260
$ws.Range("K66").Value = 1666676300
$ws.Range("M66").Value = -1666672868

$ws.Range("H102").Value = 4118639.5
$ws.Range("I102").Value = 5292679.5
$ws.Range("K102").Value = 5292679.5
$ws.Range("M102").Value = -5291057.5

$ws.Range("H122").Value = 1835721.9
$ws.Range("I122").Value = 2140842.2
$ws.Range("K122").Value = 6422526.600000001
$ws.Range("M122").Value = -6420076.600000001

$ws.Range("H135").Value = 56480.184
$ws.Range("J135").Value = 56480.184
$ws.Range("L135").Value = 56480.184
$ws.Range("N135").Value = -66620.18400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1010.73914
$ws.Range("I94").Value = 874.4666999999999
$ws.Range("K94").Value = 874.4666999999999
$ws.Range("M94").Value = -423.4666999999999

$ws.Range("H107").Value = 500603.66
$ws.Range("I107").Value = 600624.4
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 600624.4
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = -598704.4
$ws.Range("N107").Value = -4340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 68.5
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 87
$ws.Range("K7").Value = 50
$ws.Range("L7").Value = 87
$ws.Range("M7").Value = 63
$ws.Range("N7").Value = -313

$ws.Range("H62").Value = 18000
$ws.Range("I62").Value = 18000
$ws.Range("J62").Value = 18000
$ws.Range("K62").Value = 18000
$ws.Range("L62").Value = 18000
$ws.Range("M62").Value = -17376
$ws.Range("N62").Value = -19248

$ws.Range("H65").Value = 18000
$ws.Range("I65").Value = 18000
$ws.Range("J65").Value = 18000
$ws.Range("K65").Value = 90000
$ws.Range("L65").Value = 90000
$ws.Range("M65").Value = -86880
$ws.Range("N65").Value = -96240

$ws.Range("H99").Value = 13891167
$ws.Range("I99").Value = 2050
$ws.Range("J99").Value = 25002460
$ws.Range("K99").Value = 2050
$ws.Range("L99").Value = 25002460
$ws.Range("M99").Value = -552
$ws.Range("N99").Value = -25005456

$ws.Range("H126").Value = 13891167
$ws.Range("I126").Value = 2050
$ws.Range("J126").Value = 25002460
$ws.Range("K126").Value = 6150
$ws.Range("L126").Value = 75007380
$ws.Range("M126").Value = -3680
$ws.Range("N126").Value = -75012320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2267.1428
$ws.Range("J34").Value = 3300
$ws.Range("L34").Value = 9900
$ws.Range("N34").Value = -10068

$ws.Range("H39").Value = 2500
$ws.Range("J39").Value = 2500
$ws.Range("L39").Value = 7500
$ws.Range("N39").Value = -8088

$ws.Range("H55").Value = 4170
$ws.Range("J55").Value = 4170
$ws.Range("L55").Value = 12510
$ws.Range("N55").Value = -12864

$ws.Range("H122").Value = 5664.391
$ws.Range("J122").Value = 12812.111
$ws.Range("L122").Value = 115308.999
$ws.Range("N122").Value = -120208.999

$ws.Range("H131").Value = 2326609
$ws.Range("I131").Value = 7143559.5
$ws.Range("J131").Value = 1184.8276
$ws.Range("K131").Value = 21430678.5
$ws.Range("L131").Value = 3554.4828
$ws.Range("M131").Value = -21425638.5
$ws.Range("N131").Value = -13634.4828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14688.125
$ws.Range("J80").Value = 2500
$ws.Range("L80").Value = 2500
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 14688.125
$ws.Range("J83").Value = 2500
$ws.Range("L83").Value = 12500
$ws.Range("N83").Value = -22484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1945.2188
$ws.Range("I22").Value = 772.8570999999999
$ws.Range("J22").Value = 2273.48
$ws.Range("K22").Value = 772.8570999999999
$ws.Range("L22").Value = 2273.48
$ws.Range("M22").Value = -477.8570999999999
$ws.Range("N22").Value = -2863.48

$ws.Range("H27").Value = 1945.2188
$ws.Range("I27").Value = 772.8570999999999
$ws.Range("J27").Value = 2273.48
$ws.Range("K27").Value = 772.8570999999999
$ws.Range("L27").Value = 2273.48
$ws.Range("M27").Value = -665.8570999999999
$ws.Range("N27").Value = -2487.48

$ws.Range("H68").Value = 43480270
$ws.Range("I68").Value = 1866
$ws.Range("J68").Value = 90911256
$ws.Range("K68").Value = 1866
$ws.Range("L68").Value = 90911256
$ws.Range("M68").Value = -1117
$ws.Range("N68").Value = -90912754

$ws.Range("H71").Value = 43480270
$ws.Range("I71").Value = 1866
$ws.Range("J71").Value = 90911256
$ws.Range("K71").Value = 9330
$ws.Range("L71").Value = 454556280
$ws.Range("M71").Value = -5586
$ws.Range("N71").Value = -454563768

$ws.Range("H82").Value = 695443.75
$ws.Range("I82").Value = 1251637
$ws.Range("K82").Value = 1251637
$ws.Range("M82").Value = -1251276

$ws.Range("H85").Value = 695443.75
$ws.Range("I85").Value = 1251637
$ws.Range("K85").Value = 1251637
$ws.Range("M85").Value = -1250389

$ws.Range("H93").Value = 783.3333
$ws.Range("I93").Value = 783.3333
$ws.Range("K93").Value = 783.3333
$ws.Range("M93").Value = 464.6667

$ws.Range("H133").Value = 55261
$ws.Range("J133").Value = 55261
$ws.Range("L133").Value = 55261
$ws.Range("N133").Value = -60321
